$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-08-23 Saturday" "2025-08-24 Sunday"

Replace-Text "439×3=" "343×2="
Replace-Text "661×2=" "860×3="
Replace-Text "883×7=" "763×3="
Replace-Text "950×2=" "683×8="
Replace-Text "453×6=" "706×5="
Replace-Text "665×3=" "632×7="
Replace-Text "150×5=" "586×8="
Replace-Text "716×6=" "429×9="
Replace-Text "597×9=" "729×4="
Replace-Text "823×4=" "870×2="
Replace-Text "378×7=" "296×5="
Replace-Text "408×9=" "628×3="
Replace-Text "251×7=" "815×6="
Replace-Text "937×7=" "685×5="
Replace-Text "237×5=" "636×8="
Replace-Text "613×7=" "236×7="
Replace-Text "459×8=" "102×2="
Replace-Text "732×2=" "770×6="
Replace-Text "427×3=" "382×3="
Replace-Text "833×7=" "945×7="
Replace-Text "850×3=" "544×5="
Replace-Text "223×6=" "585×5="
Replace-Text "400×7=" "569×5="
Replace-Text "382×6=" "975×6="
Replace-Text "141×6=" "394×6="
